$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 49775
$ws.Range("J63").Value = 49775
$ws.Range("L63").Value = 49775
$ws.Range("N63").Value = -51023

$ws.Range("H66").Value = 49775
$ws.Range("J66").Value = 49775
$ws.Range("L66").Value = 149325
$ws.Range("N66").Value = -155565

$ws.Range("H125").Value = 1850
$ws.Range("J125").Value = 1850
$ws.Range("L125").Value = 16650
$ws.Range("N125").Value = -21570

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4441.32
$ws.Range("I32").Value = 3833.0105
$ws.Range("J32").Value = 15999.2
$ws.Range("K32").Value = 3833.0105
$ws.Range("L32").Value = 15999.2
$ws.Range("M32").Value = -3546.0105
$ws.Range("N32").Value = -16573.2

$ws.Range("H74").Value = 5503.5625
$ws.Range("I74").Value = 6671.4165
$ws.Range("K74").Value = 6671.4165
$ws.Range("M74").Value = -5797.4165

$ws.Range("H77").Value = 5503.5625
$ws.Range("I77").Value = 6671.4165
$ws.Range("K77").Value = 33357.0825
$ws.Range("M77").Value = -28989.0825

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2230.9092
$ws.Range("I99").Value = 1200
$ws.Range("J99").Value = 2334
$ws.Range("K99").Value = 1200
$ws.Range("L99").Value = 2334
$ws.Range("M99").Value = 298
$ws.Range("N99").Value = -5330

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 9175.556
$ws.Range("I4").Value = 6000
$ws.Range("J4").Value = 10082.857
$ws.Range("K4").Value = 6000
$ws.Range("L4").Value = 10082.857
$ws.Range("M4").Value = -5888
$ws.Range("N4").Value = -10306.857

$ws.Range("H31").Value = 2256
$ws.Range("I31").Value = 1364.1052
$ws.Range("J31").Value = 3424.6897
$ws.Range("K31").Value = 1364.1052
$ws.Range("L31").Value = 3424.6897
$ws.Range("M31").Value = -1069.1052
$ws.Range("N31").Value = -4014.6897

$ws.Range("H34").Value = 2256
$ws.Range("I34").Value = 1364.1052
$ws.Range("J34").Value = 3424.6897
$ws.Range("K34").Value = 1364.1052
$ws.Range("L34").Value = 3424.6897
$ws.Range("M34").Value = -1162.1052
$ws.Range("N34").Value = -3828.6897

$ws.Range("H99").Value = 10003079
$ws.Range("I99").Value = 3130
$ws.Range("K99").Value = 3130
$ws.Range("M99").Value = -1632

$ws.Range("H107").Value = 1097.9412
$ws.Range("I107").Value = 333.18182
$ws.Range("J107").Value = 2500
$ws.Range("K107").Value = 333.18182
$ws.Range("L107").Value = 2500
$ws.Range("M107").Value = 1586.81818
$ws.Range("N107").Value = -6340

$ws.Range("H126").Value = 10003079
$ws.Range("I126").Value = 3130
$ws.Range("K126").Value = 9390
$ws.Range("M126").Value = -6920

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 983.5789
$ws.Range("I131").Value = 516.4286
$ws.Range("J131").Value = 1135.6744
$ws.Range("K131").Value = 1549.2858
$ws.Range("L131").Value = 3407.023200000001
$ws.Range("M131").Value = 3490.7142
$ws.Range("N131").Value = -13487.0232

$ws.Range("H133").Value = 3632.543
$ws.Range("I133").Value = 2076
$ws.Range("J133").Value = 4255.16
$ws.Range("K133").Value = 6228
$ws.Range("L133").Value = 12765.48
$ws.Range("M133").Value = -1168
$ws.Range("N133").Value = -22885.48

$ws.Range("H140").Value = 1385.1111
$ws.Range("I140").Value = 1302.1333
$ws.Range("K140").Value = 3906.3999
$ws.Range("M140").Value = 1273.6001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 15978
$ws.Range("J109").Value = 15978
$ws.Range("L109").Value = 15978
$ws.Range("N109").Value = -18058

$ws.Range("H126").Value = 12294
$ws.Range("I126").Value = 11580
$ws.Range("K126").Value = 34740
$ws.Range("M126").Value = -32270

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3866.6667
$ws.Range("I7").Value = 3400
$ws.Range("K7").Value = 3400
$ws.Range("M7").Value = -3288

$ws.Range("H40").Value = 3115
$ws.Range("I40").Value = 2643.077
$ws.Range("K40").Value = 2643.077
$ws.Range("M40").Value = -2507.077

$ws.Range("H55").Value = 254.54546
$ws.Range("I55").Value = 171.42857
$ws.Range("J55").Value = 400
$ws.Range("K55").Value = 171.42857
$ws.Range("L55").Value = 400
$ws.Range("M55").Value = 1.571429999999992
$ws.Range("N55").Value = -746

$ws.Range("H100").Value = 4174.385
$ws.Range("I100").Value = 4009
$ws.Range("J100").Value = 4367.3335
$ws.Range("K100").Value = 4009
$ws.Range("L100").Value = 4367.3335
$ws.Range("M100").Value = -3468
$ws.Range("N100").Value = -5449.3335

$ws.Range("H126").Value = 3866.6667
$ws.Range("I126").Value = 3400
$ws.Range("K126").Value = 10200
$ws.Range("M126").Value = -7730

$ws.Range("H132").Value = 20149.621
$ws.Range("J132").Value = 4838.533
$ws.Range("L132").Value = 14515.599
$ws.Range("N132").Value = -19575.599

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3516.8572
$ws.Range("I62").Value = 2550
$ws.Range("J62").Value = 4054
$ws.Range("K62").Value = 2550
$ws.Range("L62").Value = 4054
$ws.Range("M62").Value = -1926
$ws.Range("N62").Value = -5302

$ws.Range("H65").Value = 3516.8572
$ws.Range("I65").Value = 2550
$ws.Range("J65").Value = 4054
$ws.Range("K65").Value = 12750
$ws.Range("L65").Value = 20270
$ws.Range("M65").Value = -9630
$ws.Range("N65").Value = -26510

$ws.Range("H100").Value = 796
$ws.Range("I100").Value = 309.6
$ws.Range("J100").Value = 1606.6666
$ws.Range("K100").Value = 619.2
$ws.Range("L100").Value = 3213.3332
$ws.Range("M100").Value = -78.20000000000005
$ws.Range("N100").Value = -4295.3332

$ws.Range("H113").Value = 706
$ws.Range("I113").Value = 352
$ws.Range("K113").Value = 1056
$ws.Range("M113").Value = 1114

$ws.Range("H126").Value = 6184.273
$ws.Range("I126").Value = 5558.5557
$ws.Range("J126").Value = 9000
$ws.Range("K126").Value = 16675.6671
$ws.Range("L126").Value = 27000
$ws.Range("M126").Value = -14205.6671
$ws.Range("N126").Value = -31940

$ws.Range("H127").Value = 29729
$ws.Range("J127").Value = 29729
$ws.Range("L127").Value = 29729
$ws.Range("N127").Value = -39649

$ws.Range("H136").Value = 589485.6
$ws.Range("I136").Value = 824124.4
$ws.Range("J136").Value = 2888.889
$ws.Range("K136").Value = 2472373.2
$ws.Range("L136").Value = 8666.667000000001
$ws.Range("M136").Value = -2469823.2
$ws.Range("N136").Value = -13766.667
